$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 1028, shifting existing rows 1028:1113 down to 1029:1114
$ws.Rows(1028).Insert()

# Populate the newly inserted row 1028 with its data
$ws.Range("A1028").Value = 4
$ws.Range("B1028").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C1028").Value = "Los Lagos"
$ws.Range("D1028").Value = 45106
$ws.Range("E1028").Value = 10
$ws.Range("F1028").Value = 100112020
$ws.Range("G1028").Value = "Tomate"
$ws.Range("H1028").Value = "Larga vida"
$ws.Range("I1028").Value = "Primera"
$ws.Range("J1028").Value = 400
$ws.Range("K1028").Value = 19000
$ws.Range("L1028").Value = 21000
$ws.Range("M1028").Value = 20000
$ws.Range("N1028").Value = "$/bandeja 18 kilos"
$ws.Range("O1028").Value = "Región de Arica y Parinacota"
$ws.Range("P1028").Value = 1111
$ws.Range("Q1028").Value = 18
$ws.Range("R1028").Value = "Hortaliza"
